$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-17 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-06-18 Sunday", 2) | Out-Null
$d.Content.Find.Execute("75-26=", $true, $true, $false, $false, $false, $true, 1, $false, "6+40=", 2) | Out-Null
$d.Content.Find.Execute("75-13=", $true, $true, $false, $false, $false, $true, 1, $false, "33+55=", 2) | Out-Null
$d.Content.Find.Execute("54+24=", $true, $true, $false, $false, $false, $true, 1, $false, "22+18=", 2) | Out-Null
$d.Content.Find.Execute("89-73=", $true, $true, $false, $false, $false, $true, 1, $false, "60-14=", 2) | Out-Null
$d.Content.Find.Execute("11-2=", $true, $true, $false, $false, $false, $true, 1, $false, "16+37=", 2) | Out-Null
$d.Content.Find.Execute("58+11=", $true, $true, $false, $false, $false, $true, 1, $false, "34-31=", 2) | Out-Null
$d.Content.Find.Execute("44-9=", $true, $true, $false, $false, $false, $true, 1, $false, "80-62=", 2) | Out-Null
$d.Content.Find.Execute("29-20=", $true, $true, $false, $false, $false, $true, 1, $false, "50+18=", 2) | Out-Null
$d.Content.Find.Execute("47+46=", $true, $true, $false, $false, $false, $true, 1, $false, "50+8=", 2) | Out-Null
$d.Content.Find.Execute("30+7=", $true, $true, $false, $false, $false, $true, 1, $false, "31-15=", 2) | Out-Null
$d.Content.Find.Execute("70-50=", $true, $true, $false, $false, $false, $true, 1, $false, "87+5=", 2) | Out-Null
$d.Content.Find.Execute("32+35=", $true, $true, $false, $false, $false, $true, 1, $false, "35+51=", 2) | Out-Null
$d.Content.Find.Execute("45-22=", $true, $true, $false, $false, $false, $true, 1, $false, "44-18=", 2) | Out-Null
$d.Content.Find.Execute("28+3=", $true, $true, $false, $false, $false, $true, 1, $false, "14+12=", 2) | Out-Null
$d.Content.Find.Execute("67+28=", $true, $true, $false, $false, $false, $true, 1, $false, "4+93=", 2) | Out-Null
$d.Content.Find.Execute("0+48=", $true, $true, $false, $false, $false, $true, 1, $false, "6+61=", 2) | Out-Null
$d.Content.Find.Execute("73-71=", $true, $true, $false, $false, $false, $true, 1, $false, "92-22=", 2) | Out-Null
$d.Content.Find.Execute("58+4=", $true, $true, $false, $false, $false, $true, 1, $false, "36-10=", 2) | Out-Null
$d.Content.Find.Execute("62-59=", $true, $true, $false, $false, $false, $true, 1, $false, "93-67=", 2) | Out-Null
$d.Content.Find.Execute("4+89=", $true, $true, $false, $false, $false, $true, 1, $false, "31+65=", 2) | Out-Null
$d.Content.Find.Execute("17+49=", $true, $true, $false, $false, $false, $true, 1, $false, "66-6=", 2) | Out-Null
$d.Content.Find.Execute("81-15=", $true, $true, $false, $false, $false, $true, 1, $false, "47-33=", 2) | Out-Null
$d.Content.Find.Execute("76-64=", $true, $true, $false, $false, $false, $true, 1, $false, "80-14=", 2) | Out-Null
$d.Content.Find.Execute("46+4=", $true, $true, $false, $false, $false, $true, 1, $false, "81-35=", 2) | Out-Null
$d.Content.Find.Execute("93-10=", $true, $true, $false, $false, $false, $true, 1, $false, "44-42=", 2) | Out-Null
$d.Content.Find.Execute("8+12=", $true, $true, $false, $false, $false, $true, 1, $false, "59-30=", 2) | Out-Null
$d.Content.Find.Execute("44+14=", $true, $true, $false, $false, $false, $true, 1, $false, "88-5=", 2) | Out-Null
$d.Content.Find.Execute("65+8=", $true, $true, $false, $false, $false, $true, 1, $false, "37+19=", 2) | Out-Null
$d.Content.Find.Execute("12+36=", $true, $true, $false, $false, $false, $true, 1, $false, "82-71=", 2) | Out-Null
$d.Content.Find.Execute("44-40=", $true, $true, $false, $false, $false, $true, 1, $false, "19+28=", 2) | Out-Null
$d.Content.Find.Execute("82-5=", $true, $true, $false, $false, $false, $true, 1, $false, "68-30=", 2) | Out-Null
$d.Content.Find.Execute("20-6=", $true, $true, $false, $false, $false, $true, 1, $false, "27-5=", 2) | Out-Null
$d.Content.Find.Execute("64-43=", $true, $true, $false, $false, $false, $true, 1, $false, "49-42=", 2) | Out-Null
$d.Content.Find.Execute("69+13=", $true, $true, $false, $false, $false, $true, 1, $false, "58-56=", 2) | Out-Null
$d.Content.Find.Execute("41-16=", $true, $true, $false, $false, $false, $true, 1, $false, "47+47=", 2) | Out-Null
$d.Content.Find.Execute("78+11=", $true, $true, $false, $false, $false, $true, 1, $false, "47-40=", 2) | Out-Null
$d.Content.Find.Execute("30+47=", $true, $true, $false, $false, $false, $true, 1, $false, "40+21=", 2) | Out-Null
$d.Content.Find.Execute("87-13=", $true, $true, $false, $false, $false, $true, 1, $false, "7-1=", 2) | Out-Null
$d.Content.Find.Execute("23+4=", $true, $true, $false, $false, $false, $true, 1, $false, "10+82=", 2) | Out-Null
$d.Content.Find.Execute("34-13=", $true, $true, $false, $false, $false, $true, 1, $false, "48+50=", 2) | Out-Null
$d.Content.Find.Execute("77-58=", $true, $true, $false, $false, $false, $true, 1, $false, "68-24=", 2) | Out-Null
$d.Content.Find.Execute("90-67=", $true, $true, $false, $false, $false, $true, 1, $false, "81-5=", 2) | Out-Null
$d.Content.Find.Execute("27-26=", $true, $true, $false, $false, $false, $true, 1, $false, "56+17=", 2) | Out-Null
$d.Content.Find.Execute("0+11=", $true, $true, $false, $false, $false, $true, 1, $false, "45-40=", 2) | Out-Null
$d.Content.Find.Execute("39+21=", $true, $true, $false, $false, $false, $true, 1, $false, "16+61=", 2) | Out-Null
$d.Content.Find.Execute("50+29=", $true, $true, $false, $false, $false, $true, 1, $false, "85+4=", 2) | Out-Null
$d.Content.Find.Execute("38-35=", $true, $true, $false, $false, $false, $true, 1, $false, "71-69=", 2) | Out-Null
$d.Content.Find.Execute("39+26=", $true, $true, $false, $false, $false, $true, 1, $false, "37+32=", 2) | Out-Null
$d.Content.Find.Execute("67-25=", $true, $true, $false, $false, $false, $true, 1, $false, "55-38=", 2) | Out-Null
$d.Content.Find.Execute("21-7=", $true, $true, $false, $false, $false, $true, 1, $false, "89-61=", 2) | Out-Null
$d.Content.Find.Execute("34+6=", $true, $true, $false, $false, $false, $true, 1, $false, "49-4=", 2) | Out-Null
$d.Content.Find.Execute("16-2=", $true, $true, $false, $false, $false, $true, 1, $false, "75-8=", 2) | Out-Null
$d.Content.Find.Execute("70+5=", $true, $true, $false, $false, $false, $true, 1, $false, "76-48=", 2) | Out-Null
$d.Content.Find.Execute("64+26=", $true, $true, $false, $false, $false, $true, 1, $false, "45-24=", 2) | Out-Null
$d.Content.Find.Execute("57+17=", $true, $true, $false, $false, $false, $true, 1, $false, "0+24=", 2) | Out-Null
$d.Content.Find.Execute("25+23=", $true, $true, $false, $false, $false, $true, 1, $false, "9+20=", 2) | Out-Null
$d.Content.Find.Execute("90+3=", $true, $true, $false, $false, $false, $true, 1, $false, "68-9=", 2) | Out-Null
$d.Content.Find.Execute("8+83=", $true, $true, $false, $false, $false, $true, 1, $false, "26+20=", 2) | Out-Null
$d.Content.Find.Execute("44-14=", $true, $true, $false, $false, $false, $true, 1, $false, "86-17=", 2) | Out-Null
$d.Content.Find.Execute("62-42=", $true, $true, $false, $false, $false, $true, 1, $false, "40+37=", 2) | Out-Null
$d.Content.Find.Execute("83-37=", $true, $true, $false, $false, $false, $true, 1, $false, "88-24=", 2) | Out-Null
$d.Content.Find.Execute("44-10=", $true, $true, $false, $false, $false, $true, 1, $false, "72+2=", 2) | Out-Null
$d.Content.Find.Execute("61-28=", $true, $true, $false, $false, $false, $true, 1, $false, "67-24=", 2) | Out-Null
$d.Content.Find.Execute("28+68=", $true, $true, $false, $false, $false, $true, 1, $false, "76-13=", 2) | Out-Null
$d.Content.Find.Execute("15+33=", $true, $true, $false, $false, $false, $true, 1, $false, "26-23=", 2) | Out-Null
$d.Content.Find.Execute("17+58=", $true, $true, $false, $false, $false, $true, 1, $false, "5-3=", 2) | Out-Null
$d.Content.Find.Execute("76-32=", $true, $true, $false, $false, $false, $true, 1, $false, "30+36=", 2) | Out-Null
$d.Content.Find.Execute("22+63=", $true, $true, $false, $false, $false, $true, 1, $false, "32-29=", 2) | Out-Null
$d.Content.Find.Execute("35+56=", $true, $true, $false, $false, $false, $true, 1, $false, "94-47=", 2) | Out-Null
$d.Content.Find.Execute("84-14=", $true, $true, $false, $false, $false, $true, 1, $false, "85-71=", 2) | Out-Null
$d.Content.Find.Execute("99-84=", $true, $true, $false, $false, $false, $true, 1, $false, "9+39=", 2) | Out-Null
$d.Content.Find.Execute("31+62=", $true, $true, $false, $false, $false, $true, 1, $false, "84-36=", 2) | Out-Null
$d.Content.Find.Execute("38-32=", $true, $true, $false, $false, $false, $true, 1, $false, "28+56=", 2) | Out-Null
$d.Content.Find.Execute("20+71=", $true, $true, $false, $false, $false, $true, 1, $false, "63-15=", 2) | Out-Null
$d.Content.Find.Execute("82-7=", $true, $true, $false, $false, $false, $true, 1, $false, "14+3=", 2) | Out-Null
$d.Content.Find.Execute("37-18=", $true, $true, $false, $false, $false, $true, 1, $false, "2+97=", 2) | Out-Null
$d.Content.Find.Execute("59+11=", $true, $true, $false, $false, $false, $true, 1, $false, "40-19=", 2) | Out-Null
$d.Content.Find.Execute("99-38=", $true, $true, $false, $false, $false, $true, 1, $false, "65+27=", 2) | Out-Null
$d.Content.Find.Execute("28+27=", $true, $true, $false, $false, $false, $true, 1, $false, "31+44=", 2) | Out-Null
$d.Content.Find.Execute("49+32=", $true, $true, $false, $false, $false, $true, 1, $false, "9+23=", 2) | Out-Null
$d.Content.Find.Execute("62+10=", $true, $true, $false, $false, $false, $true, 1, $false, "52+19=", 2) | Out-Null
$d.Content.Find.Execute("6+42=", $true, $true, $false, $false, $false, $true, 1, $false, "3+13=", 2) | Out-Null
$d.Content.Find.Execute("53-29=", $true, $true, $false, $false, $false, $true, 1, $false, "98-74=", 2) | Out-Null
$d.Content.Find.Execute("37-17=", $true, $true, $false, $false, $false, $true, 1, $false, "17+51=", 2) | Out-Null
$d.Content.Find.Execute("54-35=", $true, $true, $false, $false, $false, $true, 1, $false, "82-11=", 2) | Out-Null
$d.Content.Find.Execute("94-59=", $true, $true, $false, $false, $false, $true, 1, $false, "78-3=", 2) | Out-Null
$d.Content.Find.Execute("48-23=", $true, $true, $false, $false, $false, $true, 1, $false, "28+14=", 2) | Out-Null
$d.Content.Find.Execute("66+31=", $true, $true, $false, $false, $false, $true, 1, $false, "87-30=", 2) | Out-Null
$d.Content.Find.Execute("47+40=", $true, $true, $false, $false, $false, $true, 1, $false, "38+6=", 2) | Out-Null
$d.Content.Find.Execute("68+15=", $true, $true, $false, $false, $false, $true, 1, $false, "46-36=", 2) | Out-Null
$d.Content.Find.Execute("27+9=", $true, $true, $false, $false, $false, $true, 1, $false, "21+46=", 2) | Out-Null
$d.Content.Find.Execute("17+63=", $true, $true, $false, $false, $false, $true, 1, $false, "0+87=", 2) | Out-Null
$d.Content.Find.Execute("65+21=", $true, $true, $false, $false, $false, $true, 1, $false, "30+28=", 2) | Out-Null
$d.Content.Find.Execute("67-42=", $true, $true, $false, $false, $false, $true, 1, $false, "97-46=", 2) | Out-Null
$d.Content.Find.Execute("30+1=", $true, $true, $false, $false, $false, $true, 1, $false, "84-43=", 2) | Out-Null
$d.Content.Find.Execute("65+25=", $true, $true, $false, $false, $false, $true, 1, $false, "12-7=", 2) | Out-Null
$d.Content.Find.Execute("16+71=", $true, $true, $false, $false, $false, $true, 1, $false, "21+72=", 2) | Out-Null
$d.Content.Find.Execute("72-53=", $true, $true, $false, $false, $false, $true, 1, $false, "76-29=", 2) | Out-Null
$d.Content.Find.Execute("79-1=", $true, $true, $false, $false, $false, $true, 1, $false, "86-35=", 2) | Out-Null
$d.Content.Find.Execute("11+41=", $true, $true, $false, $false, $false, $true, 1, $false, "73-18=", 2) | Out-Null
